$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear B2 and B4 (they contained literal 0; after clearing, their dependent
# formulas evaluate to NULL instead of "0",)
$ws.Range("B2").ClearContents()
$ws.Range("B4").ClearContents()

# Update the view: scroll to show column L at top-left and select AM2:AM4
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("AM2:AM4").Select()

# Hide columns S through AL (19-38)
$ws.Range("S1:AL1").EntireColumn.Hidden = $true
